# NCB-530 "still not test" edit
#
# 1. Rename the "Name" lookup tab to "NAME".
# 2. On the main NCB_530 sheet, the PAYMENT_CHANNEL column (E2:E4) is
#    changed from "ALL" to the concatenation of every value on the
#    PAYMENT_CHANNEL reference sheet ("PC_1/PC_2/PC_3/PC_4/PC_5"), which
#    also widens column E to fit the longer text.
# 3. The cursor/selection left behind on several sheets while the author
#    was cross-checking values is updated to match what was captured in
#    the saved workbook.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Name" -> "NAME" ------------------------------------------
$wsName = $wb.Worksheets.Item("Name")
$wsName.Name = "NAME"

# --- 2. Update PAYMENT_CHANNEL values on NCB_530 -------------------------
$wsMain = $wb.Worksheets.Item("NCB_530")
$wsMain.Range("E2").Value = "PC_1/PC_2/PC_3/PC_4/PC_5"
$wsMain.Range("E3").Value = "PC_1/PC_2/PC_3/PC_4/PC_5"
$wsMain.Range("E4").Value = "PC_1/PC_2/PC_3/PC_4/PC_5"

# Column E grows to fit the new, longer text (~21.57 -> ~28.57 characters).
$wsMain.Columns.Item(5).ColumnWidth = 27.67

# --- 3. Restore the selections captured in each sheet ---------------------
$wsNAME = $wb.Worksheets.Item("NAME")
$wsNAME.Activate()
$wsNAME.Range("F16").Select() | Out-Null

$wsChannel = $wb.Worksheets.Item("PAYMENT_CHANNEL")
$wsChannel.Activate()
$wsChannel.Range("A2:A6").Select() | Out-Null

$wsCurrency = $wb.Worksheets.Item("MT4_CURRENCY")
$wsCurrency.Activate()
$wsCurrency.Range("I33").Select() | Out-Null

# Leave the workbook back on the main sheet, which is the one that was
# active/selected when the file was saved.
$wsMain.Activate()
$wsMain.Range("D15").Select() | Out-Null
